$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# September row data (row 10): BPJS = 89, Umum = 24, Total = SUM(B10:C10)
$ws.Range("B10").Value = 89
$ws.Range("C10").Value = 24
$ws.Range("D10").Formula = "=SUM(B10:C10)"

# Matches the author's cursor move to E13 after entering the new row
$ws.Range("E13").Select()
